$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 181; this shifts existing rows 181-330 down to 182-331
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with the new weekly data point
$ws.Cells.Item(181, 1).Value2 = 3
$ws.Cells.Item(181, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value2 = "Coquimbo"
$ws.Cells.Item(181, 4).Value2 = 44729
$ws.Cells.Item(181, 5).Value2 = 5
$ws.Cells.Item(181, 6).Value2 = 100112039
$ws.Cells.Item(181, 7).Value2 = "Ciboulette"
$ws.Cells.Item(181, 8).Value2 = "Sin especificar"
$ws.Cells.Item(181, 9).Value2 = "Primera"
$ws.Cells.Item(181, 10).Value2 = 120
$ws.Cells.Item(181, 11).Value2 = 1500
$ws.Cells.Item(181, 12).Value2 = 1500
$ws.Cells.Item(181, 13).Value2 = 1500
$ws.Cells.Item(181, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(181, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(181, 16).Value2 = 500
$ws.Cells.Item(181, 17).Value2 = 3
$ws.Cells.Item(181, 18).Value2 = "Hortaliza"
